$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.211.70"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.332.28"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'398.59"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").Value = "'125.04"
$ws.Range("E6").Value = "  +6.94%  "
$ws.Range("D7").Value = "'0.586"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.653"
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("D10").Value = "'0.118"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").Value = "'40.52"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "3.880.36"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "'8.17"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "'19.11"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "3.337.60"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "61.125.45"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'11.12"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").Value = "'0.997"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'0.0000126"
$ws.Range("E20").Value = "  +9.65%  "
$ws.Range("D21").Value = "'3.17"
$ws.Range("E21").Value = "  -6.49%  "
$ws.Range("D22").Value = "'79.54"
$ws.Range("E22").Value = "  +7.04%  "
$ws.Range("D23").Value = "'12.65"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'296.80"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'3.08"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "'4.73"
$ws.Range("E26").Value = "  +11.53%  "
$ws.Range("D27").Value = "'8.17"
$ws.Range("E27").Value = "  +8.34%  "
$ws.Range("D28").Value = "'28.77"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'0.112"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'11.19"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").Value = "'2.48"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'40.59"
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D36").Value = "'0.0472"
$ws.Range("E36").Value = "  -3.95%  "
$ws.Range("D37").Value = "'51.89"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "'3.36"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").Value = "'2.86"
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").Value = "'136.20"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'1.95"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'0.276"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").Value = "'16.50"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'3.82"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'21.02"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "2.106.97"
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "3.667.77"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = "  -3.23%  "
